$d = $word.ActiveDocument

# 1) Centro educativo: fill in the placeholder with "IES EJEMPLO"
$rng = $d.Content
$rng.Find.Execute("Centro educativo: _______________________________") | Out-Null
$rng.Text = "Centro educativo: "
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("IES EJEMPLO") | Out-Null
$rng.Bold = 1
$rng.Bold = 0

# 2) Alumno/a: fill in the placeholder with "Fulanito"
$rng = $d.Content
$rng.Find.Execute("Alumno/a: _______________________________________") | Out-Null
$rng.Text = "Alumno/a: "
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Fulanito") | Out-Null
$rng.Bold = 1
$rng.Bold = 0

# 3) Tutor/a: fill in the placeholder with "Menganita"
$rng = $d.Content
$rng.Find.Execute("Tutor/a: ________________________________________") | Out-Null
$rng.Text = "Tutor/a: "
$rng.Collapse(0) | Out-Null
$rng.InsertAfter("Menganita") | Out-Null
$rng.Bold = 1
$rng.Bold = 0

# 4) Add "Bla bla bla" to the empty paragraph right under "1. Título del proyecto"
$rng = $d.Content
$rng.Find.Execute("1. Título del proyecto") | Out-Null
$target = $rng.Start
$foundIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $target) {
        $foundIndex = $i
    }
}
$bodyPara = $d.Paragraphs.Item($foundIndex + 1)
$prng = $bodyPara.Range
$prng.Collapse(0) | Out-Null
$prng.MoveEnd(1, -1) | Out-Null
$prng.Collapse(0) | Out-Null
$prng.InsertAfter("Bla bla bla") | Out-Null
$prng.Bold = 1
$prng.Bold = 0
